$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.95
$ws.Range("I2").Value = 3.8
$ws.Range("U2").Value = 1.58
$ws.Range("AE2").Value = 12
$ws.Range("N3").Value = 19
$ws.Range("U3").Value = 1.72
$ws.Range("V4").Value = 1.67
$ws.Range("V5").Value = 1.63
$ws.Range("BD5").Value = 126
$ws.Range("N9").Value = 8
$ws.Range("Q9").Value = 2.4
$ws.Range("R9").Value = 1.53
$ws.Range("Q21").Value = 2.5
$ws.Range("R21").Value = 1.5
$ws.Range("U33").Value = 1.8
$ws.Range("V33").Value = 1.8
$ws.Range("U34").Value = 1.92
$ws.Range("V34").Value = 1.77
$ws.Range("N41").Value = 8
$ws.Range("N42").Value = 9
$ws.Range("I45").Value = 1.9
$ws.Range("N45").Value = 10
$ws.Range("O45").Value = 1.33
$ws.Range("P45").Value = 3.25
$ws.Range("R45").Value = 1.75
$ws.Range("S45").Value = 1.44
$ws.Range("T45").Value = 2.63
$ws.Range("W45").Value = 10
$ws.Range("AT45").Value = 2.63
$ws.Range("O46").Value = 1.3
$ws.Range("P46").Value = 3.4
$ws.Range("Q46").Value = 2.03
$ws.Range("R46").Value = 1.83
$ws.Range("I47").Value = 1.73
$ws.Range("M48").Value = 1.02
$ws.Range("O48").Value = 1.11
$ws.Range("M49").Value = 1.03
$ws.Range("O49").Value = 1.22
$ws.Range("U50").Value = 1.77
$ws.Range("V50").Value = 1.92
$ws.Range("V51").Value = 1.72
$ws.Range("M52").Value = 1.05
$ws.Range("O52").Value = 1.3
$ws.Range("M53").Value = 1.05
$ws.Range("O53").Value = 1.33
$ws.Range("M54").Value = 1.05
$ws.Range("O54").Value = 1.3
$ws.Range("M55").Value = 1.03
$ws.Range("O55").Value = 1.17
$ws.Range("M56").Value = 1.03
$ws.Range("O56").Value = 1.19
$ws.Range("M57").Value = 1.01
$ws.Range("O57").Value = 1.1
$ws.Range("Q57").Value = 1.44
$ws.Range("R57").Value = 2.7
$ws.Range("G58").Value = 2.25
$ws.Range("I58").Value = 3.5
$ws.Range("J58").Value = 2.88
$ws.Range("N58").Value = 9
$ws.Range("X58").Value = 10
$ws.Range("AD58").Value = 6
$ws.Range("AJ58").Value = 13
$ws.Range("AK58").Value = 41
$ws.Range("AO58").Value = 12
$ws.Range("G61").Value = 2.01
$ws.Range("M61").Value = 1.02
$ws.Range("N61").Value = 21
$ws.Range("G62").Value = 1.53
$ws.Range("G63").Value = 2.05
$ws.Range("I63").Value = 3.3
$ws.Range("AD63").Value = 7
$ws.Range("AQ63").Value = 34
$ws.Range("G64").Value = 1.76
$ws.Range("O64").Value = 1.18
$ws.Range("P64").Value = 4.5
$ws.Range("Q64").Value = 1.6
$ws.Range("R64").Value = 2.3
$ws.Range("G65").Value = 2.7
$ws.Range("I65").Value = 2.45
$ws.Range("J65").Value = 3.25
$ws.Range("L65").Value = 3
$ws.Range("U65").Value = 1.57
$ws.Range("V65").Value = 2.25
$ws.Range("Z65").Value = 29
$ws.Range("AA65").Value = 21
$ws.Range("AK65").Value = 23
$ws.Range("AN65").Value = 5
$ws.Range("AO65").Value = 15
$ws.Range("Q69").Value = 1.63
$ws.Range("U69").Value = 1.58
$ws.Range("Q70").Value = 1.54
$ws.Range("U70").Value = 1.47
$ws.Range("R71").Value = 1.62
$ws.Range("U71").Value = 1.8
$ws.Range("V71").Value = 1.8
$ws.Range("R72").Value = 1.63
$ws.Range("R73").Value = 1.63
$ws.Range("Q74").Value = 1.69
$ws.Range("R74").Value = 2.07
$ws.Range("R75").Value = 1.44
$ws.Range("Q76").Value = 1.77
$ws.Range("R77").Value = 1.67
$ws.Range("M80").Value = 1.05
$ws.Range("O80").Value = 1.37
$ws.Range("U80").Value = 1.8
$ws.Range("V80").Value = 1.8
$ws.Range("M82").Value = 1.04
$ws.Range("N82").Value = 9
$ws.Range("O83").Value = 1.06
$ws.Range("P83").Value = 8
$ws.Range("Q84").Value = 1.6
$ws.Range("R84").Value = 2.3
$ws.Range("M87").Value = 1.07
$ws.Range("N87").Value = 9
